# Updates cryptos price/volume cells to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cell to source a clean "General"/default style (s=0) from,
# used to strip the Text-format style that Excel attaches to cells we
# force to stay text (values that otherwise parse as a plain number).
$blankFormatSource = $ws.Range("Z1")

$ws.Range("D2").Value = "66.528.59"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.576.00"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.01"
$blankFormatSource.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.84"
$blankFormatSource.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$blankFormatSource.Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").Value = "2.576.83"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$blankFormatSource.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$blankFormatSource.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.16"
$blankFormatSource.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.81"
$blankFormatSource.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "3.058.84"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$blankFormatSource.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "66.603.90"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "2.593.37"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.47"
$blankFormatSource.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  -4.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.69"
$blankFormatSource.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  -4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.70"
$blankFormatSource.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$blankFormatSource.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.57"
$blankFormatSource.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.17"
$blankFormatSource.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  -7.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$blankFormatSource.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "68.95"
$blankFormatSource.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").Value = "2.733.19"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$blankFormatSource.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "0.0₃0981"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.65"
$blankFormatSource.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$blankFormatSource.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  -3.33%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$blankFormatSource.Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$blankFormatSource.Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$blankFormatSource.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  -4.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.84"
$blankFormatSource.Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.71"
$blankFormatSource.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.360"
$blankFormatSource.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.23"
$blankFormatSource.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.78"
$blankFormatSource.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$blankFormatSource.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$blankFormatSource.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("D46").Value = "0.0₆0288"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.76"
$blankFormatSource.Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.565"
$blankFormatSource.Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$blankFormatSource.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$blankFormatSource.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0761"
$blankFormatSource.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  -1.74%  "

$excel.CutCopyMode = $false
